# Daily attendance processing - 2025-12-20 03:08:59
# Normalize the "Recorded By" (column G) lists so that the "System" entry
# is moved one position earlier (swapped with the entry immediately
# preceding it), for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $text = [string]$val
    if ($text -notlike "*System*") {
        continue
    }

    $parts = $text -split ", "
    $idx = [array]::IndexOf($parts, "System")

    if ($idx -gt 0) {
        $tmp = $parts[$idx - 1]
        $parts[$idx - 1] = $parts[$idx]
        $parts[$idx] = $tmp
        $cell.Value2 = $parts -join ", "
    }
}
